$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row7 B: 110uF Cap link -> update to new 0603 cap URL
$ws.Range("B7").Value = "https://www.lcsc.com/product-detail/Multilayer-Ceramic-Capacitors-MLCC-SMD-SMT_CCTC-TCC0603X7R104K500CT_C282519.html"
$ws.Hyperlinks.Item(1).Address
